{"js": "// The document previously held the tag \"<id>\" as one run, the literal\n// text \"p051r_2\" as a separately-formatted run, and the closing tag\n// \"</id>\" as another run matching the first run's (Courier New /\n// color 7f6000) formatting. The edit collapses all three runs into a\n// single run containing \"<id>p051r_2</id>\", picking up the formatting\n// of the run that originally held the surrounding tags.\n\nconst body = context.document.body;\n\n// Word's search matches across run boundaries against the rendered\n// text, so this single search locates the three runs as one logical\n// hit even though \"<id>\", \"p051r_2\" and \"</id>\" live in separate\n// <w:r> elements in the underlying OOXML.\nconst results = body.search(\"<id>p051r_2</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p051r_2</id>' in the document body.\");\n}\n\n// Replacing the whole matched range with the same visible text merges\n// the underlying runs into one, inheriting the formatting of the\n// first run in the match (the Courier-New tag styling), exactly like\n// typing the replacement text over the selection in Word.\nconst target = results.items[0];\ntarget.insertText(\"<id>p051r_2</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document previously held the tag \"<id>\" as one run, the literal\n# text \"p051r_2\" as a separately-formatted run, and the closing tag\n# \"</id>\" as another run matching the first run's (Courier New /\n# color 7f6000) formatting. The edit collapses all three runs into a\n# single run containing \"<id>p051r_2</id>\", picking up the formatting\n# of the run that originally held the surrounding tags.\n\n$d = $word.ActiveDocument\n\n# Word's Find matches across run boundaries against the rendered text,\n# so this single Find/Replace locates the three runs as one logical\n# hit even though \"<id>\", \"p051r_2\" and \"</id>\" live in separate runs\n# in the underlying OOXML, and rewrites that whole span in one go.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n# Replace = 2 -> wdReplaceAll: rewrite every match with ReplaceWith,\n# which merges the matched runs into one run carrying the formatting\n# of the first run in the match (the Courier-New tag styling) -\n# exactly like using Word's Find & Replace dialog.\n$result = $rng.Find.Execute(\"<id>p051r_2</id>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<id>p051r_2</id>\", 2)\n"}
